$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 73: Thu Aug 15 2024, T, 6 hours, notes ---
$a73 = $ws.Cells.Item(73, 1)
$a73.Value = 45519
$a73.NumberFormat = $ws.Cells.Item(72, 1).NumberFormat

$b73 = $ws.Cells.Item(73, 2)
$b73.Value = "T"
$b73.NumberFormat = $ws.Cells.Item(72, 2).NumberFormat

$ws.Cells.Item(73, 3).Value = 6

$e73 = $ws.Cells.Item(73, 5)
$e73.Value = "Editting the rmd file and knitting. Corrections include: adding links to data and models, adding graphs for each model, cleaning up and displaying all dataset formats, reducing output overall to only include relevant information."
$e73.HorizontalAlignment = -4108
$e73.WrapText = $true

# --- Row 74: Fri Aug 16 2024, F, 8 hours, link text ---
$a74 = $ws.Cells.Item(74, 1)
$a74.Value = 45520
$a74.NumberFormat = $ws.Cells.Item(72, 1).NumberFormat

$b74 = $ws.Cells.Item(74, 2)
$b74.Value = "F"
$b74.NumberFormat = $ws.Cells.Item(72, 2).NumberFormat

$ws.Cells.Item(74, 3).Value = 8

$g74 = $ws.Cells.Item(74, 7)
$g74.Value = "text for results, organize rdata files to folder, add a description to the github"

# --- Merge the notes cell across the two rows ---
$ws.Range("E73:E74").Merge()

# --- Update selection to match the saved view ---
[void]$ws.Range("G74").Select()

Write-Host "done"
